$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix the four "Arrow Functions" code-block strings in row 8: the original text was
# missing the opening parenthesis before "str1.length)" in the console.log(...) call.
$ws.Range("E8").Value = 'const myFunction = () => {\n const str1 = \''some text\''\n console.log(str1.length)\n}'
$ws.Range("F8").Value = 'const myFunction = __ __ {\n _____ str1 = \''some text\''\n console.___(str1.length)\n}'
$ws.Range("G8").Value = 'var myFunction = (is perfect) = {\n str1 = some text\n console.slog(str1.lenth)\n}'
$ws.Range("H8").Value = 'const myFunction = () => {Æ const str1 = \''some text\''Æ console.log(str1.length)Æ}'

$ws.Range("H8").Select() | Out-Null
